# Update @base / @prefix values on the "Info" sheet to require angle
# brackets around URIs (D1, D2 and D3 need to be wrapped in <...>).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Range("D1").Value = "<http://sales.data/purchases/2015>"
$ws.Range("D2").Value = "<http://sales.data/purchases#>"
$ws.Range("D3").Value = "<http://sales.data/schema#>"

$ws.Range("D4").Select()
